# Update EL-NumberValidation "NumberError" sheet: mark all test rows as
# executed by filling in Result ("Pass") and Date (execution timestamp)
# for rows 2-61, columns A and B respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$dates = @(
    "Mon Apr 08 17:45:08 EDT 2024",
    "Mon Apr 08 17:45:17 EDT 2024",
    "Mon Apr 08 17:45:27 EDT 2024",
    "Mon Apr 08 17:45:36 EDT 2024",
    "Mon Apr 08 17:45:45 EDT 2024",
    "Mon Apr 08 17:45:55 EDT 2024",
    "Mon Apr 08 17:46:04 EDT 2024",
    "Mon Apr 08 17:46:13 EDT 2024",
    "Mon Apr 08 17:46:23 EDT 2024",
    "Mon Apr 08 17:46:32 EDT 2024",
    "Mon Apr 08 17:46:41 EDT 2024",
    "Mon Apr 08 17:46:50 EDT 2024",
    "Mon Apr 08 17:46:59 EDT 2024",
    "Mon Apr 08 17:47:08 EDT 2024",
    "Mon Apr 08 17:47:18 EDT 2024",
    "Mon Apr 08 17:47:27 EDT 2024",
    "Mon Apr 08 17:47:36 EDT 2024",
    "Mon Apr 08 17:47:46 EDT 2024",
    "Mon Apr 08 17:47:55 EDT 2024",
    "Mon Apr 08 17:48:04 EDT 2024",
    "Mon Apr 08 17:48:13 EDT 2024",
    "Mon Apr 08 17:48:22 EDT 2024",
    "Mon Apr 08 17:48:31 EDT 2024",
    "Mon Apr 08 17:48:40 EDT 2024",
    "Mon Apr 08 17:48:49 EDT 2024",
    "Mon Apr 08 17:48:59 EDT 2024",
    "Mon Apr 08 17:49:08 EDT 2024",
    "Mon Apr 08 17:49:17 EDT 2024",
    "Mon Apr 08 17:49:27 EDT 2024",
    "Mon Apr 08 17:49:36 EDT 2024",
    "Mon Apr 08 17:49:45 EDT 2024",
    "Mon Apr 08 17:49:55 EDT 2024",
    "Mon Apr 08 17:50:04 EDT 2024",
    "Mon Apr 08 17:50:13 EDT 2024",
    "Mon Apr 08 17:50:22 EDT 2024",
    "Mon Apr 08 17:50:32 EDT 2024",
    "Mon Apr 08 17:50:42 EDT 2024",
    "Mon Apr 08 17:50:51 EDT 2024",
    "Mon Apr 08 17:51:02 EDT 2024",
    "Mon Apr 08 17:51:11 EDT 2024",
    "Mon Apr 08 17:51:20 EDT 2024",
    "Mon Apr 08 17:51:29 EDT 2024",
    "Mon Apr 08 17:51:38 EDT 2024",
    "Mon Apr 08 17:51:47 EDT 2024",
    "Mon Apr 08 17:51:56 EDT 2024",
    "Mon Apr 08 17:52:07 EDT 2024",
    "Mon Apr 08 17:52:16 EDT 2024",
    "Mon Apr 08 17:52:25 EDT 2024",
    "Mon Apr 08 17:52:35 EDT 2024",
    "Mon Apr 08 17:52:44 EDT 2024",
    "Mon Apr 08 17:52:53 EDT 2024",
    "Mon Apr 08 17:53:03 EDT 2024",
    "Mon Apr 08 17:53:12 EDT 2024",
    "Mon Apr 08 17:53:21 EDT 2024",
    "Mon Apr 08 17:53:30 EDT 2024",
    "Mon Apr 08 17:53:39 EDT 2024",
    "Mon Apr 08 17:53:48 EDT 2024",
    "Mon Apr 08 17:53:57 EDT 2024",
    "Mon Apr 08 17:54:06 EDT 2024",
    "Mon Apr 08 17:54:15 EDT 2024"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $resultCell = $ws.Cells.Item($row, 1)
    $resultCell.Value = "Pass"
    $resultCell.Style = "Normal"
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.Value = $dates[$i]
    $dateCell.Style = "Normal"
}
